# queuing_theory.xlsx - apply diff changes across Part 1..Part 5
# All cells in this workbook are stored as text (inlineStr) even when the
# content looks numeric, so every numeric-looking value is written with a
# leading apostrophe to force Excel to keep it as Text instead of
# auto-converting it to a Number cell.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Part 1
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Part 1")

$ws1.Range("B2").Value = "'0.001"
$ws1.Range("D2").Value = "'12"
$ws1.Range("B3").Value = "'20.0"
$ws1.Range("D3").Value = "'0.0006416879652056608"

# C4 / D4 were already-empty inlineStr cells; the new sheet drops them
# entirely.
$ws1.Range("C4").ClearContents()
$ws1.Range("D4").ClearContents()

# ----------------------------------------------------------------------
# Part 2 - gains a Results/Values block in columns C:D (dimension A1:B9 -> A1:D9)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Part 2")

$ws2.Range("C1").Value = "Results"
$ws2.Range("D1").Value = "Values"

$ws2.Range("A2").Value = "Max P(wait)"
$ws2.Range("B2").Value = "'0.9"
$ws2.Range("C2").Value = "Number of Servers"
$ws2.Range("D2").Value = "'2"

$ws2.Range("A3").Value = "Max E(w)"
$ws2.Range("B3").Value = "'30.0"
$ws2.Range("C3").Value = "E(S)"
$ws2.Range("D3").Value = "'0.02380952380952381"

$ws2.Range("A4").Value = "Arrival Rate"
$ws2.Range("B4").Value = "'15.0"
$ws2.Range("C4").Value = "E(N)"
$ws2.Range("D4").Value = "'0.02380952380952381"

$ws2.Range("A5").Value = "Service Rate"
$ws2.Range("B5").Value = "'21.0"

# Rows 6-9 (Results / Number of Servers / E(S) / E(N)) are untouched.

# ----------------------------------------------------------------------
# Part 3
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Part 3")

$ws3.Range("B2").Value = "'2.0"
$ws3.Range("D2").Value = "'0.6701030927835052"
$ws3.Range("B3").Value = "'2.0"
$ws3.Range("D3").Value = "'1.8556701030927834"
$ws3.Range("B4").Value = "'0.8"
$ws3.Range("D4").Value = "'0.538888888888889"
$ws3.Range("D5").Value = "'0.32638888888888884"

# ----------------------------------------------------------------------
# Part 4
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Part 4")

$ws4.Range("B2").Value = "'10.0"
$ws4.Range("D2").Value = "'0.6780485778822252"
$ws4.Range("B3").Value = "'4.0"
$ws4.Range("D3").Value = "'1.5243927764721845"
$ws4.Range("B4").Value = "'5.0"

# ----------------------------------------------------------------------
# Part 5
# ----------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Part 5")

$ws5.Range("B2").Value = "'10.0"
$ws5.Range("D2").Value = "'1.051948051948052"
$ws5.Range("B3").Value = "'5.0"
$ws5.Range("D3").Value = "'-0.051948051948051965"
$ws5.Range("B4").Value = "'4.0"
